$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9042444229125977
$ws.Range("B1").Value = 1.28155255317688
$ws.Range("C1").Value = 3.522947788238525
$ws.Range("D1").Value = 3.406658411026001
$ws.Range("E1").Value = 0.4981203377246857
